$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 49 (row 114 block) - Maandag row: Carlo (col D) hours corrected 0 -> 4
$ws.Range("D115").Value = 4

# Week 50 (row 122 block) - Maandag row: Carlo (col D) hours corrected 4 -> 0
$ws.Range("D123").Value = 0

# Week 51 (row 130 block) - newly filled-in hours for Dinsdag / Woensdag / Donderdag rows
$ws.Range("B132").Value = 5
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 5
$ws.Range("E132").Value = 5
$ws.Range("F132").Value = 5
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 5
$ws.Range("I132").Value = 5

$ws.Range("B133").Value = 3
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 3
$ws.Range("E133").Value = 2
$ws.Range("F133").Value = 3
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 3
$ws.Range("I133").Value = 0

$ws.Range("B134").Value = 2
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 2
$ws.Range("E134").Value = 2
$ws.Range("F134").Value = 2
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = 2
$ws.Range("I134").Value = 2

# Week 50 block - individual lesuren note, L141 filled in
$ws.Range("L141").Value = 2

# Update the active selection to match the authored edit location
$ws.Range("D124").Select()
